# Update the two-digit / one-digit division practice problems in the table.
#
# Each problem cell is addressed by its (row, column) position rather than a
# document-wide Find/Replace, because some problem text (e.g. "42÷3=") occurs
# more than once in the table and a global Find/Replace cannot target a single
# occurrence reliably. Assigning Range.Text directly replaces only that cell's
# contents (Word keeps the trailing cell-mark) and preserves the run formatting
# (font/size) already present in the cell.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$replacements = @(
    @{Row=1; Col=1; Old="65÷4="; New="55÷6="}
    @{Row=1; Col=2; Old="98÷3="; New="91÷6="}
    @{Row=1; Col=3; Old="73÷8="; New="42÷7="}
    @{Row=1; Col=4; Old="95÷2="; New="16÷5="}
    @{Row=1; Col=5; Old="51÷5="; New="70÷9="}
    @{Row=5; Col=1; Old="42÷3="; New="92÷4="}
    @{Row=5; Col=2; Old="54÷2="; New="51÷9="}
    @{Row=5; Col=3; Old="41÷9="; New="32÷7="}
    @{Row=5; Col=4; Old="24÷9="; New="38÷2="}
    @{Row=5; Col=5; Old="14÷9="; New="66÷8="}
    @{Row=9; Col=1; Old="93÷4="; New="15÷4="}
    @{Row=9; Col=2; Old="62÷8="; New="11÷2="}
    @{Row=9; Col=3; Old="71÷5="; New="86÷9="}
    @{Row=9; Col=4; Old="68÷5="; New="78÷5="}
    @{Row=9; Col=5; Old="80÷3="; New="57÷8="}
    @{Row=13; Col=1; Old="19÷6="; New="72÷4="}
    @{Row=13; Col=2; Old="54÷8="; New="10÷5="}
    @{Row=13; Col=3; Old="55÷7="; New="62÷2="}
    @{Row=13; Col=4; Old="38÷6="; New="97÷4="}
    @{Row=13; Col=5; Old="99÷9="; New="71÷7="}
    @{Row=17; Col=1; Old="91÷2="; New="16÷3="}
    @{Row=17; Col=2; Old="33÷9="; New="36÷2="}
    @{Row=17; Col=3; Old="42÷3="; New="27÷3="}
    @{Row=17; Col=4; Old="60÷3="; New="85÷9="}
    @{Row=17; Col=5; Old="68÷3="; New="31÷6="}
)

foreach ($rep in $replacements) {
    $cellRange = $t.Cell($rep.Row, $rep.Col).Range
    # Cell().Range.Text includes the trailing cell-mark (CR + BEL); strip it
    # before comparing against the expected plain problem text.
    $current = $cellRange.Text.TrimEnd([char]13, [char]7)
    if ($current -eq $rep.Old) {
        $cellRange.Text = $rep.New
    } else {
        Write-Output ("Unexpected text at R" + $rep.Row + "C" + $rep.Col + ": " + $current)
    }
}

Write-Output "Done."
